$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grow the Table1 ListObject by 5 rows so the table range / autofilter
# expand the way Excel does when new rows are added inside a table.
$lo = $ws.ListObjects.Item(1)
for ($i = 0; $i -lt 5; $i++) {
    $lo.ListRows.Add() | Out-Null
}

# Rewrite rows 52-63: four brand-new rows for the board-mount hardware
# (76-79), the pre-existing 90-96 rows shifted down two rows, the old
# "Wire Anchor" (96) renamed/refiled as the T-slot version, and a new
# row 97 for the XY joint wire guide revision.
$ws.Range("A52").Value() = 76
$ws.Range("B52").Value() = "Electrical"
$ws.Range("C52").Value() = "Mount"
$ws.Range("D52").Value() = "L"
$ws.Range("E52").Value() = "Control Board Clamp"
$ws.Range("F52").Value() = "ABS"
$ws.Range("G52").Value() = 1
$ws.Range("I52").Value() = "76 - Electrical - Mount - L Control Board Clamp.stl"
$ws.Range("A53").Value() = 77
$ws.Range("B53").Value() = "Electrical"
$ws.Range("C53").Value() = "Mount"
$ws.Range("D53").Value() = "R"
$ws.Range("E53").Value() = "Control Board Clamp"
$ws.Range("F53").Value() = "ABS"
$ws.Range("G53").Value() = 1
$ws.Range("I53").Value() = "77 - Electrical - Mount - R Control Board Clamp.stl"
$ws.Range("A54").Value() = 78
$ws.Range("B54").Value() = "Electrical"
$ws.Range("C54").Value() = "Board Tray"
$ws.Range("D54").Value() = "N"
$ws.Range("E54").Value() = "Octopus Board Tray"
$ws.Range("F54").Value() = "ABS"
$ws.Range("G54").Value() = 1
$ws.Range("I54").Value() = "78 - Electrical - Board Tray - Octopus Board Tray.stl"
$ws.Range("A55").Value() = 79
$ws.Range("B55").Value() = "Electrical"
$ws.Range("C55").Value() = "Board Tray"
$ws.Range("D55").Value() = "N"
$ws.Range("E55").Value() = "SKR Board Tray"
$ws.Range("F55").Value() = "ABS"
$ws.Range("G55").Value() = 1
$ws.Range("I55").Value() = "79 - Electrical - Board Tray - SKR Board Tray.stl"
$ws.Range("A56").Value() = 90
$ws.Range("B56").Value() = "Misc"
$ws.Range("C56").Value() = "Frame"
$ws.Range("D56").Value() = "N"
$ws.Range("E56").Value() = "Bracket Cover"
$ws.Range("F56").Value() = "ABS"
$ws.Range("G56").Value() = "?"
$ws.Range("I56").Value() = "90 - Misc - Bracket Cover.stl"
$ws.Range("A57").Value() = 91
$ws.Range("B57").Value() = "Misc"
$ws.Range("C57").Value() = "Motion"
$ws.Range("D57").Value() = "N"
$ws.Range("E57").Value() = "Limit Switch Cap"
$ws.Range("F57").Value() = "ABS"
$ws.Range("G57").Value() = 5
$ws.Range("I57").Value() = "91 - Misc - Limit Switch Cap.stl"
$ws.Range("A58").Value() = 92
$ws.Range("B58").Value() = "Misc"
$ws.Range("C58").Value() = "Frame"
$ws.Range("D58").Value() = "N"
$ws.Range("E58").Value() = "Frame Base Foot"
$ws.Range("F58").Value() = "TPU"
$ws.Range("G58").Value() = 4
$ws.Range("I58").Value() = "92 - Misc - Frame Base Foot (TPU).stl"
$ws.Range("A59").Value() = 93
$ws.Range("B59").Value() = "XY"
$ws.Range("C59").Value() = "Motion"
$ws.Range("D59").Value() = "R"
$ws.Range("E59").Value() = "X Limit Switch Mount"
$ws.Range("F59").Value() = "ABS"
$ws.Range("G59").Value() = 1
$ws.Range("I59").Value() = "93 - XY - Motion - X Limit Switch Mount.stl"
$ws.Range("A60").Value() = 94
$ws.Range("B60").Value() = "Misc"
$ws.Range("C60").Value() = "Tools"
$ws.Range("D60").Value() = "N"
$ws.Range("E60").Value() = "Gantry Lock"
$ws.Range("F60").Value() = "ABS"
$ws.Range("G60").Value() = 2
$ws.Range("I60").Value() = "94 - Misc - Tools - Gantry Lock.stl"
$ws.Range("A61").Value() = 95
$ws.Range("B61").Value() = "Misc"
$ws.Range("C61").Value() = "Wiring"
$ws.Range("D61").Value() = "N"
$ws.Range("E61").Value() = "Wire Guide"
$ws.Range("F61").Value() = "TPU"
$ws.Range("G61").Value() = "?"
$ws.Range("H61").Value() = "Still Testing"
$ws.Range("I61").Value() = "95 - Misc - Wiring - Wire Guide.stl"
$ws.Range("A62").Value() = 96
$ws.Range("B62").Value() = "Misc"
$ws.Range("C62").Value() = "Wiring"
$ws.Range("D62").Value() = "N"
$ws.Range("E62").Value() = "T Slot Wire Anchor"
$ws.Range("F62").Value() = "ABS"
$ws.Range("G62").Value() = 10
$ws.Range("I62").Value() = "96 - Misc - Wiring - T Slot Wire Anchor.stl"
$ws.Range("A63").Value() = 97
$ws.Range("B63").Value() = "Misc"
$ws.Range("C63").Value() = "Wiring"
$ws.Range("D63").Value() = "R"
$ws.Range("E63").Value() = "XY Joint Wire Guide"
$ws.Range("F63").Value() = "ABS"
$ws.Range("G63").Value() = 1
$ws.Range("I63").Value() = "97 - Misc - Wiring - R XY Joint Wire Guide.stl"

# Reflect the cursor/selection shown in the saved file.
$ws.Range("H62").Select()
